$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Period" value cell (B4) to pass client timezone using
# from/to .toString() formatting instead of the old String.format call.
$ws.Range("B4").Value = '${from.toString("YYYY.MM.dd HH:mm:ss")+" - "+to.toString("YYYY.MM.dd HH:mm:ss")}'
